$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 13
$ws.Cells.Item(2, 2).Value = "Aircraft in Attitude mode . Unable to hover . Fly with caution ."
$ws.Cells.Item(2, 3).Value = "Unable to hover"
$ws.Cells.Item(2, 4).Value = "5-7"
$ws.Cells.Item(2, 5).Value = "Missing"

$ws.Cells.Item(3, 1).Value = 15
$ws.Cells.Item(3, 2).Value = "Aircraft is close to the Home Point . Initiating Return to Home will now trigger Auto Landing ."
$ws.Cells.Item(3, 3).Value = "Initiating Return to Home will now trigger Auto Landing"
$ws.Cells.Item(3, 4).Value = "8-16"
$ws.Cells.Item(3, 5).Value = "Missing"

$ws.Cells.Item(4, 1).Value = 29
$ws.Cells.Item(4, 2).Value = "Another aircraft is dangerously close , please descend to a safer altitude ."
$ws.Cells.Item(4, 3).Value = "Another aircraft is dangerously close"
$ws.Cells.Item(4, 4).Value = "0-4"
$ws.Cells.Item(4, 5).Value = "Missing"

$ws.Cells.Item(5, 1).Value = 48
$ws.Cells.Item(5, 2).Value = "Check whether propellers are installed correctly . If the propellers are installed correctly and the aircraft still cannot takeoff, a motor error may exist . Contact DJI Support for assistance ."
$ws.Cells.Item(5, 3).Value = "If the propellers are installed correctly and the aircraft still cannot takeoff, a motor error may exist"
$ws.Cells.Item(5, 4).Value = "7-23"
$ws.Cells.Item(5, 5).Value = "Missing"

$ws.Cells.Item(6, 1).Value = 48
$ws.Cells.Item(6, 2).Value = "Check whether propellers are installed correctly . If the propellers are installed correctly and the aircraft still cannot takeoff, a motor error may exist . Contact DJI Support for assistance ."
$ws.Cells.Item(6, 3).Value = "If the propellers are installed correctly and the aircraft still cannot takeoff,"
$ws.Cells.Item(6, 4).Value = "7-18"
$ws.Cells.Item(6, 5).Value = "'False"

$ws.Cells.Item(7, 1).Value = 50
$ws.Cells.Item(7, 2).Value = "Compass abnormal . Solution: 1. Ensure there are no magnets or metal objects near the aircraft . The ground or walls may contain metal . Move away from sources of interference before attempting flight . 2. Calibrate Compass Before Takeoff ."
$ws.Cells.Item(7, 3).Value = "2. Calibrate Compass Before Takeoff"
$ws.Cells.Item(7, 4).Value = "35-39"
$ws.Cells.Item(7, 5).Value = "Missing"

$ws.Cells.Item(8, 1).Value = 50
$ws.Cells.Item(8, 2).Value = "Compass abnormal . Solution: 1. Ensure there are no magnets or metal objects near the aircraft . The ground or walls may contain metal . Move away from sources of interference before attempting flight . 2. Calibrate Compass Before Takeoff ."
$ws.Cells.Item(8, 3).Value = "Calibrate Compass Before Takeoff"
$ws.Cells.Item(8, 4).Value = "36-39"
$ws.Cells.Item(8, 5).Value = "'False"

$ws.Cells.Item(9, 1).Value = 70
$ws.Cells.Item(9, 2).Value = "Downward ambient light too low . Obstacle avoidance unavailable . Fly with caution . Backward ambient light too low . Backward obstacle avoidance unavailable . Only infrared sensors available . Fly with caution ."
$ws.Cells.Item(9, 3).Value = "Backward ambient light too low"
$ws.Cells.Item(9, 4).Value = "14-18"
$ws.Cells.Item(9, 5).Value = "Missing"

$ws.Cells.Item(10, 1).Value = 77
$ws.Cells.Item(10, 2).Value = "Exiting GPS mode : Unknown Error ."
$ws.Cells.Item(10, 3).Value = "Exiting GPS mode"
$ws.Cells.Item(10, 4).Value = "0-2"
$ws.Cells.Item(10, 5).Value = "Missing"

$ws.Cells.Item(11, 1).Value = 81
$ws.Cells.Item(11, 2).Value = "Extra payload detected . Return aircraft to an area nearby the home point promptly and fly in a wind-free environment to ensure flight safety ."
$ws.Cells.Item(11, 3).Value = "Return aircraft to an area nearby the home point promptly and fly in a wind-free environment to ensure flight safety"
$ws.Cells.Item(11, 4).Value = "4-23"
$ws.Cells.Item(11, 5).Value = "Missing"

$ws.Cells.Item(12, 1).Value = 91
$ws.Cells.Item(12, 2).Value = "GEO Zone Info: The target area is in an Altitude Zone . Flight altitude restricted to nnn ."
$ws.Cells.Item(12, 3).Value = "GEO Zone Info: The target area is in an Altitude Zone"
$ws.Cells.Item(12, 4).Value = "0-10"
$ws.Cells.Item(12, 5).Value = "Missing"

$ws.Cells.Item(13, 1).Value = 91
$ws.Cells.Item(13, 2).Value = "GEO Zone Info: The target area is in an Altitude Zone . Flight altitude restricted to nnn ."
$ws.Cells.Item(13, 3).Value = "GEO Zone Info:"
$ws.Cells.Item(13, 4).Value = "0-2"
$ws.Cells.Item(13, 5).Value = "Missing"

$ws.Cells.Item(14, 1).Value = 91
$ws.Cells.Item(14, 2).Value = "GEO Zone Info: The target area is in an Altitude Zone . Flight altitude restricted to nnn ."
$ws.Cells.Item(14, 3).Value = "The target area is in an Altitude Zone"
$ws.Cells.Item(14, 4).Value = "3-10"
$ws.Cells.Item(14, 5).Value = "'False"

$ws.Cells.Item(15, 1).Value = 92
$ws.Cells.Item(15, 2).Value = "GEO: You are in a Warning Zone (Airport Class Airspace Unpaved Airports Power Plant) . Fly with caution ."
$ws.Cells.Item(15, 3).Value = "GEO: You are in a Warning Zone (Airport Class Airspace Unpaved Airports Power Plant)"
$ws.Cells.Item(15, 4).Value = "0-13"
$ws.Cells.Item(15, 5).Value = "Missing"

$ws.Cells.Item(16, 1).Value = 115
$ws.Cells.Item(16, 2).Value = "Landin ."
$ws.Cells.Item(16, 3).Value = "Landin"
$ws.Cells.Item(16, 4).Value = "0-0"
$ws.Cells.Item(16, 5).Value = "Missing"

$ws.Cells.Item(17, 1).Value = 126
$ws.Cells.Item(17, 2).Value = "No GPS . Return to Home failed . Return the aircraft to home manually ."
$ws.Cells.Item(17, 3).Value = "Return to Home failed"
$ws.Cells.Item(17, 4).Value = "3-6"
$ws.Cells.Item(17, 5).Value = "Missing"

$ws.Cells.Item(18, 1).Value = 131
$ws.Cells.Item(18, 2).Value = "Panorama Captured Successfully ."
$ws.Cells.Item(18, 3).Value = "Panorama Captured Successfully"
$ws.Cells.Item(18, 4).Value = "0-2"
$ws.Cells.Item(18, 5).Value = "Missing"

$ws.Cells.Item(19, 1).Value = 131
$ws.Cells.Item(19, 2).Value = "Panorama Captured Successfully ."
$ws.Cells.Item(19, 3).Value = "Panorama Captured"
$ws.Cells.Item(19, 4).Value = "0-1"
$ws.Cells.Item(19, 5).Value = "'False"
